$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.227.74'
$ws.Range('E2').Value = '  +0.38%  '
$ws.Range('D3').Value = '1.897.03'
$ws.Range('E3').Value = '  -0.40%  '
$ws.Range('E4').Value = '  +0.39%  '
$ws.Range('D5').Value = '''307.53'
$ws.Range('E5').Value = '  +0.02%  '
$ws.Range('E6').Value = '  +0.27%  '
$ws.Range('D7').Value = '''0.5204'
$ws.Range('E7').Value = '  -0.25%  '
$ws.Range('D8').Value = '''0.3763'
$ws.Range('E8').Value = '  -0.46%  '
$ws.Range('D9').Value = '''0.07289'
$ws.Range('E9').Value = '  +0.90%  '
$ws.Range('E10').Value = '  -0.16%  '
$ws.Range('D11').Value = '''0.9006'
$ws.Range('E11').Value = '  +0.78%  '
$ws.Range('D12').Value = '''0.08184'
$ws.Range('E12').Value = '  +6.67%  '
$ws.Range('D13').Value = '''96.66'
$ws.Range('E13').Value = '  +2.30%  '
$ws.Range('D14').Value = '1.891.68'
$ws.Range('E14').Value = '  -0.66%  '
$ws.Range('D15').Value = '''5.281'
$ws.Range('E15').Value = '  +0.80%  '
$ws.Range('D16').Value = '''1.003'
$ws.Range('E16').Value = '  +0.47%  '
$ws.Range('D17').Value = '''0.000008609'
$ws.Range('E17').Value = '  +0.97%  '
$ws.Range('D18').Value = '''14.54'
$ws.Range('E18').Value = '  +0.17%  '
$ws.Range('D19').Value = '''1.002'
$ws.Range('E19').Value = '  +0.30%  '
$ws.Range('D20').Value = '27.259.62'
$ws.Range('E20').Value = '  +0.30%  '
$ws.Range('D21').Value = '''5.088'
$ws.Range('E21').Value = '  +0.30%  '
$ws.Range('D22').Value = '''10.71'
$ws.Range('E22').Value = '  +0.84%  '
$ws.Range('D23').Value = '''6.406'
$ws.Range('E23').Value = '  -0.28%  '
$ws.Range('D24').Value = '''147.71'
$ws.Range('E24').Value = '  +1.31%  '
$ws.Range('E25').Value = '  -0.23%  '
$ws.Range('E26').Value = '  +0.67%  '
$ws.Range('E27').Value = '  +0.74%  '
$ws.Range('D28').Value = '''115.24'
$ws.Range('E28').Value = '  +0.50%  '
$ws.Range('D29').Value = '''4.929'
$ws.Range('E29').Value = '  -0.72%  '
$ws.Range('D30').Value = '''4.827'
$ws.Range('E30').Value = '  +0.37%  '
$ws.Range('D31').Value = '''0.09222'
$ws.Range('E31').Value = '  +0.16%  '
$ws.Range('E32').Value = '  -0.56%  '
$ws.Range('D33').Value = '''0.7947'
$ws.Range('E33').Value = '  +2.20%  '
$ws.Range('D34').Value = '''1.220'
$ws.Range('E34').Value = '  -1.79%  '
$ws.Range('D35').Value = '''3.453'
$ws.Range('E35').Value = '  +4.65%  '
$ws.Range('D36').Value = '''2.951'
$ws.Range('E36').Value = '  -1.07%  '
$ws.Range('D37').Value = '''2.590'
$ws.Range('E37').Value = '  -0.71%  '
$ws.Range('E38').Value = '  +0.11%  '
$ws.Range('E39').Value = '  -0.45%  '
$ws.Range('D40').Value = '''1.072'
$ws.Range('E40').Value = '  -0.11%  '
$ws.Range('D41').Value = '''8.952'
$ws.Range('E41').Value = '  -0.68%  '
$ws.Range('D42').Value = '''6.562'
$ws.Range('E42').Value = '  -1.11%  '
$ws.Range('D43').Value = '''115.32'
$ws.Range('E43').Value = '  -3.16%  '
$ws.Range('D44').Value = '''0.1515'
$ws.Range('E44').Value = '  -0.30%  '
$ws.Range('D45').Value = '''0.4899'
$ws.Range('E45').Value = '  +1.17%  '
$ws.Range('D46').Value = '''1.002'
$ws.Range('E46').Value = '  +0.28%  '
$ws.Range('D47').Value = '''10.03'
$ws.Range('E47').Value = '  -1.80%  '
$ws.Range('D48').Value = '''1.624'
$ws.Range('E48').Value = '  +1.27%  '
$ws.Range('D49').Value = '''38.25'
$ws.Range('E49').Value = '  +1.83%  '
$ws.Range('E50').Value = '  -1.15%  '
$ws.Range('D51').Value = '''0.05943'
$ws.Range('E51').Value = '  +0.38%  '
